# Fix mojibake: replace "Â±" (UTF-8 bytes for U+00B1 mis-decoded as
# Windows-1252/Latin-1, i.e. U+00C2 U+00B1) with the correct
# "±" (U+00B1 PLUS-MINUS SIGN) in the f1_score_weighted, training_time
# and test_time columns (B2:D17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$badChar = "$([char]0xC2)$([char]0xB1)"
$goodChar = "$([char]0xB1)"

$range = $ws.Range("B2:D17")

foreach ($cell in $range.Cells) {
    $value = $cell.Value()
    if ($null -ne $value -and $value -is [string] -and $value.Contains($badChar)) {
        $cell.Value = $value.Replace($badChar, $goodChar)
    }
}
